$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos table on Sheet1 (columns B:Coin, C:Link, D:Price,
# E:Volume(1h)) with the latest scraped prices/24h changes, matching the
# automated "Updated cryptos list ... with GitHub Actions" commit.
# Two coin rows (13/14 and 35/36) also swapped rank order.

# Cells whose new text values would otherwise be auto-parsed by Excel as
# numbers (e.g. "237.95"). Force text storage via NumberFormat, then revert
# the cell style back to Normal so no extra style index sticks around.
$textCells = @("D5", "D8", "D11", "D14", "D16", "D18", "D20", "D21", "D25", "D26", "D27", "D31", "D33", "D35", "D36", "D37", "D40", "D42", "D43", "D44", "D47")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = '237.95'
$ws.Range("D8").Value = '41.82'
$ws.Range("D11").Value = '0.0988'
$ws.Range("D14").Value = '11.38'
$ws.Range("D16").Value = '4.70'
$ws.Range("D18").Value = '69.95'
$ws.Range("D20").Value = '240.54'
$ws.Range("D21").Value = '12.20'
$ws.Range("D25").Value = '169.04'
$ws.Range("D26").Value = '7.97'
$ws.Range("D27").Value = '1.81'
$ws.Range("D31").Value = '0.0555'
$ws.Range("D33").Value = '3.99'
$ws.Range("D35").Value = '0.822'
$ws.Range("D36").Value = '2.00'
$ws.Range("D37").Value = '1.31'
$ws.Range("D40").Value = '89.74'
$ws.Range("D42").Value = '14.86'
$ws.Range("D43").Value = '12.87'
$ws.Range("D44").Value = '2.31'
$ws.Range("D47").Value = '0.0551'

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining cell updates (names, links, percentages, and non-ambiguous text).
$ws.Range("D2").Value = '35.134.34'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.853.33'
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("E4").Value = '  +0.56%  '
$ws.Range("E5").Value = '  +2.97%  '
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("E8").Value = '  +4.33%  '
$ws.Range("E9").Value = '  +1.14%  '
$ws.Range("E10").Value = '  +1.23%  '
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("D12").Value = '2.122.20'
$ws.Range("E12").Value = '  +1.81%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.886.66'
$ws.Range("E13").Value = '  +3.75%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("E15").Value = '  +1.03%  '
$ws.Range("E16").Value = '  +1.14%  '
$ws.Range("D17").Value = '35.121.65'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("E25").Value = '  -2.90%  '
$ws.Range("E26").Value = '  +2.14%  '
$ws.Range("E27").Value = '  +18.21%  '
$ws.Range("E28").Value = '  +1.32%  '
$ws.Range("E29").Value = '  -0.62%  '
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("E31").Value = '  +0.68%  '
$ws.Range("E32").Value = '  -0.78%  '
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("E34").Value = '  +27.88%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("E35").Value = '  +17.56%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("E36").Value = '  +9.58%  '
$ws.Range("E37").Value = '  +5.22%  '
$ws.Range("E38").Value = '  +8.10%  '
$ws.Range("E39").Value = '  +3.08%  '
$ws.Range("E40").Value = '  -3.65%  '
$ws.Range("D41").Value = '1.339.02'
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("E42").Value = '  +0.95%  '
$ws.Range("E43").Value = '  +53.04%  '
$ws.Range("E44").Value = '  +1.58%  '
$ws.Range("E45").Value = '  +0.40%  '
$ws.Range("E46").Value = '  -0.58%  '
$ws.Range("E47").Value = '  +5.92%  '
$ws.Range("E48").Value = '  +2.99%  '
$ws.Range("D49").Value = '2.039.79'
$ws.Range("E49").Value = '  +1.94%  '
$ws.Range("E50").Value = '  +1.27%  '
$ws.Range("E51").Value = '  +0.47%  '
